$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '20.530.70'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.474.70'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.45%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9550'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '278.07'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3618'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3055'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.42'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.057'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06642'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.512'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.186'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9538'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001029'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.474.73'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.16%  '
$ws.Range('E19').Value = '  +5.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.492'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.48'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.12'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.265'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '20.571.65'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.94'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.130'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.17'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.634.98'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '113.54'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.952'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.022'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8093'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07985'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.510'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.225'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05863'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.733'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02049'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.91%  '
$ws.Range('B40').Value = 'Frax'
$ws.Range('C40').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9534'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.71%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.34'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1878'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.457'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5293'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.20%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.26'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.92%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.526'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '118.25'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5189'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.814'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06471'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9784'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.66%  '
